$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "278.02"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "6.42%"
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "27.21"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "0.70%"
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "4.865"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "3.79%"
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.06249"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "0.52%"
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "6.859"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "1.60%"
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.8795"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "3.24%"
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.9446"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "3.03%"
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.1450"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "3.36%"
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.05166"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "6.92%"
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07345"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "3.69%"
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.03155"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "1.65%"
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.09042"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "-0.14%"
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.001558"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "1.84%"
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0006284"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "2.22%"
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.005888"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "-1.87%"
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.453"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "0.27%"
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "2.81%"
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "2.266"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "4.66%"
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "-0.62%"
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.852"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "-5.85%"
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04310"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "1.91%"
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.001176"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "-2.36%"
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.004271"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "4.84%"
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.0001201"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "-0.01%"
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "3.07%"
$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.04019"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "0.74%"
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.006707"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "62.57%"
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1151"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "3.40%"
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.002161"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "-2.27%"
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.01235"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "-7.10%"
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00005089"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "-1.41%"
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00000000750"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "-0.06%"
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.372"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "890.60%"
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.00002100"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "-0.06%"
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0002000"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "-0.06%"
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = "6"

$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = "6"
